$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new rows to the terminology table
# (order of writes matches the order new shared-string entries were created)
$ws.Range("B15").Value = "Read"
$ws.Range("B16").Value = "Write"
$ws.Range("C15").Value = "From, Load"
$ws.Range("C16").Value = "To, Save"
$ws.Range("A15").Value = "System serialization Read"
$ws.Range("A16").Value = "System serialization Write"

# Update selection to match the committed state
$ws.Range("A16").Select()
